$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RES installed")
$ws.Range("C3").Value = 10
